$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the cells in the same order they were originally typed so that
# the shared-string table indices come out in the expected sequence.
$ws.Range("A1").Value = "sadfs"
$ws.Range("C3").Value = "dsf"
$ws.Range("D4").Value = "sdfa"
$ws.Range("E5").Value = "dsfa"
$ws.Range("F6").Value = "asdff"
$ws.Range("G6").Value = "a"
$ws.Range("G7").Value = "dfa"
$ws.Range("H7").Value = "f"
$ws.Range("H8").Value = "asd"
$ws.Range("I8").Value = "f"
$ws.Range("I9").Value = "f"
$ws.Range("I10").Value = "adfasd"
$ws.Range("I11").Value = "f"
$ws.Range("H11").Value = "sd"
$ws.Range("G11").Value = "f"
$ws.Range("F11").Value = "sdf"
$ws.Range("E11").Value = "a"
$ws.Range("D11").Value = "ffass"
$ws.Range("B12").Value = "f"
$ws.Range("C12").Value = "dsf"
$ws.Range("A13").Value = "f"
$ws.Range("B13").Value = "as"
$ws.Range("A14").Value = "as"

# Leave the selection on the last-entered cell, matching the saved view state.
[void]$ws.Range("A14").Select()
